$wb = $excel.ActiveWorkbook

# --- Sheet: normality ---
$ws = $wb.Worksheets.Item("normality")
$ws.Range("D3").Value = 0.9624
$ws.Range("E3").Value = 0.8173
$ws.Range("D4").Value = 0.9314
$ws.Range("E4").Value = 0.3953
$ws.Range("D5").Value = 0.8181
$ws.Range("E5").Value = 0.0152
$ws.Range("D6").Value = 0.9167999999999999
$ws.Range("E6").Value = 0.2606
$ws.Range("D7").Value = 0.9609
$ws.Range("E7").Value = 0.7959000000000001
$ws.Range("D8").Value = 0.9631
$ws.Range("E8").Value = 0.8267
$ws.Range("D9").Value = 0.9
$ws.Range("E9").Value = 0.1339
$ws.Range("D10").Value = 0.949
$ws.Range("E10").Value = 0.5839
$ws.Range("D11").Value = 0.9473
$ws.Range("E11").Value = 0.5587
$ws.Range("D12").Value = 0.9529
$ws.Range("E12").Value = 0.6433
$ws.Range("D13").Value = 0.9126
$ws.Range("E13").Value = 0.1987
$ws.Range("D14").Value = 0.8824
$ws.Range("E14").Value = 0.0769

# --- Sheet: pairwise_tests ---
$ws = $wb.Worksheets.Item("pairwise_tests")

# Numeric columns I, L, N, O
$ws.Range("I3").Value = -1.9064
$ws.Range("L3").Value = 0.083
$ws.Range("N3").Value = -0.487
$ws.Range("O3").Value = 0.249

$ws.Range("I4").Value = -5.1551
$ws.Range("L4").Value = 0.0003
$ws.Range("N4").Value = -1.5555
$ws.Range("O4").Value = 0.0009

$ws.Range("I5").Value = -5.5744
$ws.Range("L5").Value = 0.0002
$ws.Range("N5").Value = -2.0139
$ws.Range("O5").Value = 0.0006000000000000001

$ws.Range("I6").Value = -2.4594
$ws.Range("L6").Value = 0.0301
$ws.Range("N6").Value = -0.6484
$ws.Range("O6").Value = 0.09029999999999999

$ws.Range("I7").Value = -5.6743
$ws.Range("L7").Value = 0.0001
$ws.Range("N7").Value = -1.7078
$ws.Range("O7").Value = 0.0003

$ws.Range("I8").Value = -6.2311
$ws.Range("L8").Value = 0
$ws.Range("N8").Value = -2.1274
$ws.Range("O8").Value = 0

# Text column M (BF10) - values are stored as text (shared strings), not numbers.
# Temporarily force text format so the numeric-looking strings aren't
# auto-converted to numbers, then restore the original (default) style
# so we don't leave a formatting change behind.
$ws.Range("M3:M8").NumberFormat = "@"
$ws.Range("M3").Value = "1.132"
$ws.Range("M4").Value = "107.004"
$ws.Range("M5").Value = "185.507"
$ws.Range("M6").Value = "2.396"
$ws.Range("M7").Value = "280.442"
$ws.Range("M8").Value = "593.948"
$ws.Range("M3:M8").Style = "Normal"
